# Remove the obsolete "collect_region_properties.lutFile" parameter row
# from the structural parameter properties sheet.
#
# That row currently lives at row 47 (A47 = "collect_region_properties.lutFile").
# Deleting the entire row shifts every subsequent row up by one, which is
# exactly what the target diff shows (row 48 -> 47, ... row 83 removed,
# dimension A1:H83 -> A1:H82, shared-string indices shift down by one for
# every string whose original index was > 64, and the now-unused
# "collect_region_properties.lutFile" shared string itself disappears).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Delete()

# Match the author's final selection/view state in the diff
# (topLeftCell scroll anchor is gone, selection moved to B38).
$ws.Range("B38").Select()
